$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The "Periodo Mora" column (E) listed the debt years in ascending order
# (2003..2007). The previous periods were removed and replaced with the
# new ones, now listed in descending order (2007..2003), i.e. the
# underlying database (shared strings) for the year values is rewritten.
$ws.Range("E16").Value = "2007"
$ws.Range("E17").Value = "2006"
$ws.Range("E18").Value = "2005"
$ws.Range("E19").Value = "2004"
$ws.Range("E20").Value = "2003"
